{"js": "// Office.js (Word JavaScript API) edit script.\n// Splits several run(s) of plain text into multiple runs interleaved with\n// <w:proofErr> spell/grammar-check markers (mirroring what Word itself\n// writes once its proofing pass has looked at freshly (re)typed text),\n// and appends a trailing \"->->\" paragraph at the end of the document.\n//\n// Body of: async (context) => { ... }\n\nconst OOXML_NS =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" ' +\n  'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>{BODY}</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>';\n\nfunction wrapParagraph(innerRunsXml, pPrXml) {\n  const pPr = pPrXml || '<w:pPr><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr>';\n  return OOXML_NS.replace('{BODY}', '<w:p>' + pPr + innerRunsXml + '</w:p>');\n}\n\n// A plain (non-list) run properties block used throughout this document.\nconst RPR = '<w:rPr><w:lang w:val=\"en-US\"/></w:rPr>';\n\nfunction run(text, preserve) {\n  const space = preserve ? ' xml:space=\"preserve\"' : '';\n  return '<w:r>' + RPR + '<w:t' + space + '>' + text + '</w:t></w:r>';\n}\n\nfunction proofErr(type) {\n  return '<w:proofErr w:type=\"' + type + '\"/>';\n}\n\nconst LIST_PPR =\n  '<w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/>' +\n  '<w:numId w:val=\"1\"/></w:numPr><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr>';\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// --- Paragraph: \"Mention bug blitz : say u could understand ... issues\"\nconst pBlitz = items.filter(p => p.text.indexOf(\"Mention bug blitz\") === 0)[0];\nconst blitzRuns =\n  run(\"Mention bug \", true) +\n  proofErr(\"gramStart\") +\n  run(\"blitz :\", false) +\n  proofErr(\"gramEnd\") +\n  run(\" say u could understand the flow of code base which helped me solve many issues in \", true) +\n  proofErr(\"spellStart\") +\n  run(\"there\", false) +\n  proofErr(\"spellEnd\") +\n  run(\" tickets were pending so every quarter  last week use to focus on bugs and issues\", true);\npBlitz.insertOoxml(wrapParagraph(blitzRuns), \"Replace\");\nawait context.sync();\n\n// --- Paragraph: \"Mention about refdata controller and how flow cmes to refdata service through api\"\nconst pRefdata = items.filter(p => p.text.indexOf(\"Mention about refdata controller\") === 0)[0];\nconst refdataRuns =\n  run(\"Mention about \", true) +\n  proofErr(\"spellStart\") +\n  run(\"refdata\", false) +\n  proofErr(\"spellEnd\") +\n  run(\" controller and how flow \", true) +\n  proofErr(\"spellStart\") +\n  run(\"cmes\", false) +\n  proofErr(\"spellEnd\") +\n  run(\" to \", true) +\n  proofErr(\"spellStart\") +\n  run(\"refdata\", false) +\n  proofErr(\"spellEnd\") +\n  run(\" service through \", true) +\n  proofErr(\"spellStart\") +\n  run(\"api\", false) +\n  proofErr(\"spellEnd\");\npRefdata.insertOoxml(wrapParagraph(refdataRuns, LIST_PPR), \"Replace\");\nawait context.sync();\n\n// --- Paragraph: \"With api endpoint search and navigate \"\nconst pWithApi = items.filter(p => p.text.indexOf(\"With api endpoint search and navigate\") === 0)[0];\nconst withApiRuns =\n  run(\"With \", true) +\n  proofErr(\"spellStart\") +\n  run(\"api\", false) +\n  proofErr(\"spellEnd\") +\n  run(\" endpoint search and navigate \", true);\npWithApi.insertOoxml(wrapParagraph(withApiRuns, LIST_PPR), \"Replace\");\nawait context.sync();\n\n// --- Paragraph: \"Mention about Wellsfargo issue\"\nconst pWells = items.filter(p => p.text.indexOf(\"Mention about Wellsfargo issue\") === 0)[0];\nconst wellsRuns =\n  run(\"Mention about \", true) +\n  proofErr(\"spellStart\") +\n  run(\"Wellsfargo\", false) +\n  proofErr(\"spellEnd\") +\n  run(\" issue\", true);\npWells.insertOoxml(wrapParagraph(wellsRuns, LIST_PPR), \"Replace\");\nawait context.sync();\n\n// --- Paragraph: \"The db division\"\nconst pDb = items.filter(p => p.text.indexOf(\"The db division\") === 0)[0];\nconst dbRuns =\n  run(\"The \", true) +\n  proofErr(\"spellStart\") +\n  run(\"db\", false) +\n  proofErr(\"spellEnd\") +\n  run(\" division\", true);\npDb.insertOoxml(wrapParagraph(dbRuns, LIST_PPR), \"Replace\");\nawait context.sync();\n\n// --- Paragraph: \"Apple s some division they wont compromise on security or any other thing...\"\nconst pApple = items.filter(p => p.text.indexOf(\"Apple s some division\") === 0)[0];\nconst appleRuns =\n  run(\"Apple s some division they \", true) +\n  proofErr(\"spellStart\") +\n  run(\"wont\", false) +\n  proofErr(\"spellEnd\") +\n  run(\" compromise on security\", true) +\n  run(\" or any other thing it was quite difficult for our team to satisfy them\", true);\npApple.insertOoxml(wrapParagraph(appleRuns, LIST_PPR), \"Replace\");\nawait context.sync();\n\n// --- New trailing paragraph \"->->\"\nconst pAppleRange = pApple.getRange(\"End\");\nconst arrowRuns = run(\"->->\", false);\npAppleRange.insertOoxml(wrapParagraph(arrowRuns), \"After\");\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Splits several run(s) of plain text into multiple runs interleaved with\n# <w:proofErr> spell/grammar-check markers (mirroring what Word itself\n# writes once its proofing pass has looked at freshly (re)typed text),\n# and appends a trailing \"->->\" paragraph at the end of the document.\n\n$d = $word.ActiveDocument\n\n$W = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"'\n$RPR = '<w:rPr><w:lang w:val=\"en-US\"/></w:rPr>'\n$LISTPPR = '<w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr>'\n$PLAINPPR = '<w:pPr><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr>'\n\nfunction Make-Run([string]$text, [bool]$preserve) {\n    if ($preserve) {\n        return '<w:r>' + $RPR + '<w:t xml:space=\"preserve\">' + $text + '</w:t></w:r>'\n    } else {\n        return '<w:r>' + $RPR + '<w:t>' + $text + '</w:t></w:r>'\n    }\n}\n\nfunction Make-ProofErr([string]$type) {\n    return '<w:proofErr w:type=\"' + $type + '\"/>'\n}\n\nfunction Make-Paragraph([string]$pPr, [string]$innerRuns) {\n    return '<w:p ' + $W + '>' + $pPr + $innerRuns + '</w:p>'\n}\n\nfunction Find-Paragraph([string]$startsWith) {\n    foreach ($p in $d.Paragraphs) {\n        if ($p.Range.Text.StartsWith($startsWith)) {\n            return $p\n        }\n    }\n    return $null\n}\n\n# --- Paragraph: \"Mention bug blitz : say u could understand ... issues\"\n$pBlitz = Find-Paragraph \"Mention bug blitz\"\n$blitzRuns = (Make-Run \"Mention bug \" $true) + (Make-ProofErr \"gramStart\") + (Make-Run \"blitz :\" $false) + (Make-ProofErr \"gramEnd\") + (Make-Run \" say u could understand the flow of code base which helped me solve many issues in \" $true) + (Make-ProofErr \"spellStart\") + (Make-Run \"there\" $false) + (Make-ProofErr \"spellEnd\") + (Make-Run \" tickets were pending so every quarter  last week use to focus on bugs and issues\" $true)\n$pBlitz.Range.InsertXML((Make-Paragraph $PLAINPPR $blitzRuns))\n\n# --- Paragraph: \"Mention about refdata controller and how flow cmes to refdata service through api\"\n$pRefdata = Find-Paragraph \"Mention about refdata controller\"\n$refdataRuns = (Make-Run \"Mention about \" $true) + (Make-ProofErr \"spellStart\") + (Make-Run \"refdata\" $false) + (Make-ProofErr \"spellEnd\") + (Make-Run \" controller and how flow \" $true) + (Make-ProofErr \"spellStart\") + (Make-Run \"cmes\" $false) + (Make-ProofErr \"spellEnd\") + (Make-Run \" to \" $true) + (Make-ProofErr \"spellStart\") + (Make-Run \"refdata\" $false) + (Make-ProofErr \"spellEnd\") + (Make-Run \" service through \" $true) + (Make-ProofErr \"spellStart\") + (Make-Run \"api\" $false) + (Make-ProofErr \"spellEnd\")\n$pRefdata.Range.InsertXML((Make-Paragraph $LISTPPR $refdataRuns))\n\n# --- Paragraph: \"With api endpoint search and navigate \"\n$pWithApi = Find-Paragraph \"With api endpoint search and navigate\"\n$withApiRuns = (Make-Run \"With \" $true) + (Make-ProofErr \"spellStart\") + (Make-Run \"api\" $false) + (Make-ProofErr \"spellEnd\") + (Make-Run \" endpoint search and navigate \" $true)\n$pWithApi.Range.InsertXML((Make-Paragraph $LISTPPR $withApiRuns))\n\n# --- Paragraph: \"Mention about Wellsfargo issue\"\n$pWells = Find-Paragraph \"Mention about Wellsfargo issue\"\n$wellsRuns = (Make-Run \"Mention about \" $true) + (Make-ProofErr \"spellStart\") + (Make-Run \"Wellsfargo\" $false) + (Make-ProofErr \"spellEnd\") + (Make-Run \" issue\" $true)\n$pWells.Range.InsertXML((Make-Paragraph $LISTPPR $wellsRuns))\n\n# --- Paragraph: \"The db division\"\n$pDb = Find-Paragraph \"The db division\"\n$dbRuns = (Make-Run \"The \" $true) + (Make-ProofErr \"spellStart\") + (Make-Run \"db\" $false) + (Make-ProofErr \"spellEnd\") + (Make-Run \" division\" $true)\n$pDb.Range.InsertXML((Make-Paragraph $LISTPPR $dbRuns))\n\n# --- Paragraph: \"Apple s some division they wont compromise on security or any other thing...\"\n$pApple = Find-Paragraph \"Apple s some division\"\n$appleRuns = (Make-Run \"Apple s some division they \" $true) + (Make-ProofErr \"spellStart\") + (Make-Run \"wont\" $false) + (Make-ProofErr \"spellEnd\") + (Make-Run \" compromise on security\" $true) + (Make-Run \" or any other thing it was quite difficult for our team to satisfy them\" $true)\n$pApple.Range.InsertXML((Make-Paragraph $LISTPPR $appleRuns))\n\n# --- New trailing paragraph \"->->\"\n$pApple = Find-Paragraph \"Apple s some division\"\n$pApple.Range.InsertParagraphAfter()\n$newP = $d.Paragraphs($d.Paragraphs.Count)\n$arrowRuns = Make-Run \"->->\" $false\n$newP.Range.InsertXML((Make-Paragraph $PLAINPPR $arrowRuns))\n"}
